$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'265.35"
$ws.Range("D4").Value = "'6.207"
$ws.Range("D5").Value = "'0.06158"
$ws.Range("D6").Value = "'3.566"
$ws.Range("D7").Value = "'6.714"
$ws.Range("D8").Value = "'1.353"
$ws.Range("D9").Value = "'0.8271"
$ws.Range("D11").Value = "'0.1598"
$ws.Range("D12").Value = "'0.08224"
$ws.Range("D13").Value = "'0.03403"
$ws.Range("D15").Value = "'0.09233"
$ws.Range("D16").Value = "'3.899"
$ws.Range("D17").Value = "'0.001714"
$ws.Range("D18").Value = "'0.04808"
$ws.Range("D19").Value = "'0.006227"
$ws.Range("D20").Value = "'0.006297"
$ws.Range("D21").Value = "'0.001099"
$ws.Range("D22").Value = "'0.0001501"
$ws.Range("D23").Value = "'3.742"
$ws.Range("D25").Value = "'0.3349"
$ws.Range("D27").Value = "'0.0002685"
$ws.Range("D40").Value = "'0.04610"
$ws.Range("D41").Value = "'0.006985"
$ws.Range("D42").Value = "'0.1137"
$ws.Range("D43").Value = "'0.003132"
$ws.Range("D44").Value = "'0.01092"
$ws.Range("D45").Value = "'0.00006170"
$ws.Range("D47").Value = "'0.7709"
$ws.Range("D48").Value = "'0.2056"
$ws.Range("D49").Value = "'0.00001502"
$ws.Range("D50").Value = "'0.01242"
